$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

$rng = $ws.Range("A2:F6")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "7.35 - y_1"
$ws.Range("B2").Value = "-7.35"
$ws.Range("D2").Value = "0.03"
$ws.Range("E2").Value = "8.2"
$ws.Range("F2").Value = "1.2"

$ws.Range("A3").Value = "-7.35 + y_1"
$ws.Range("B3").Value = "3.3499999999999996"
$ws.Range("D3").Value = "0.85"
$ws.Range("E3").Value = "3.5999999999999996"
$ws.Range("F3").Value = "4.0"

$ws.Range("A4").Value = "-3.7499999999999964 - 2x + y_1 + 4y_2"
$ws.Range("B4").Value = "-12.250000000000004"
$ws.Range("D4").Value = "0.73"
$ws.Range("E4").Value = "5.8"
$ws.Range("F4").Value = "9.0"

$ws.Range("A5").Value = "-67.23 + 8x + y_1"
$ws.Range("B5").Value = "18.550000000000004"
$ws.Range("D5").Value = "0.02"
$ws.Range("E5").Value = "6.7"
$ws.Range("F5").Value = "5.4"

$ws.Range("A6").Value = "-5.5 - 2x - 2y_1"
$ws.Range("B6").Value = "-17.5"
$ws.Range("D6").Value = "0.72"
$ws.Range("E6").Value = "8.0"
$ws.Range("F6").Value = "5.300000000000001"

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")

$rng = $ws.Range("A2:C2")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "7.4"
$ws.Range("B2").Value = "7.35"
$ws.Range("C2").Value = "2.8"

# --- Sheet: Vector_bf ---
$ws = $wb.Worksheets.Item("Vector_bf")

$rng = $ws.Range("A2:A3")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "0.8700000000000001"
$ws.Range("A3").Value = "-2.92"

# --- Sheet: Vector_BF (index 6; name lookup is case-insensitive and would
#     collide with "Vector_bf") ---
$ws = $wb.Worksheets.Item(6)

$rng = $ws.Range("A2:A4")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "-25.0"
$ws.Range("A3").Value = "11.1"
$ws.Range("A4").Value = "-25.2"
